# "AI conversion and adjustment" -- extend the ln(x+1) table on Sheet1
# with two more duration samples (23000 and 10000000) and update the
# sheet's view state (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New data rows 22 and 23, continuing the B/C pattern used by rows 17-21
# (B = raw value, C = LN(B+1)).
$ws.Range("B22").Value = 23000
$ws.Range("C22").Formula = "=LN(B22+1)"

$ws.Range("B23").Value = 10000000
$ws.Range("C23").Formula = "=LN(B23+1)"

# Recalculate so the new formulas carry cached results.
$excel.Calculate()

# Scroll the view down a bit and move the selection to C10.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C10").Select()
